$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 673 ("「我らも米を食います」..."), which shifts all following rows up by one.
$ws.Rows.Item(673).Delete()
